# Append 9 more "master-reg_center_user" rows (10002-10010 / 110021-110029)
# below the existing data (which already ends at row 21), mirroring the
# pattern used for the first block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow  = 22
$firstRegId   = 10002
$firstUsrId   = 110021
$rowCount     = 9

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstNewRow + $i
    $ws.Cells.Item($r, 1).Value = $firstRegId + $i
    $ws.Cells.Item($r, 2).Value = $firstUsrId + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Match the post-edit selection: the user clicked the row under the new
# data (row 31) which, with the whole row selected, spans to the end of
# the sheet.
[void]$ws.Range("A31:XFD1048576").Select()

# Page setup was switched to Portrait orientation with a print job queued
# at 300 dpi.
$ws.PageSetup.Orientation = 1
